$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 4000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 4000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 4000
$ws.Range("N13").Value = -4338
$ws.Range("M13").ClearContents()
$ws.Range("H32").Value = 10002
$ws.Range("J32").Value = 10002
$ws.Range("L32").Value = 10002
$ws.Range("N32").Value = -10654
$ws.Range("H33").Value = 10870889
$ws.Range("I33").Value = 15625326
$ws.Range("K33").Value = 15625326
$ws.Range("M33").Value = -15625097
$ws.Range("H62").Value = 3687
$ws.Range("I62").Value = 3410.9167
$ws.Range("K62").Value = 3410.9167
$ws.Range("M62").Value = -2786.9167
$ws.Range("H65").Value = 3687
$ws.Range("I65").Value = 3410.9167
$ws.Range("K65").Value = 17054.5835
$ws.Range("M65").Value = -13934.5835
$ws.Range("H82").Value = 2137.375
$ws.Range("I82").Value = 1014.1429
$ws.Range("J82").Value = 10000
$ws.Range("K82").Value = 3042.4287
$ws.Range("L82").Value = 30000
$ws.Range("M82").Value = -2636.4287
$ws.Range("N82").Value = -30812
$ws.Range("H85").Value = 2137.375
$ws.Range("I85").Value = 1014.1429
$ws.Range("J85").Value = 10000
$ws.Range("K85").Value = 3042.4287
$ws.Range("L85").Value = 30000
$ws.Range("M85").Value = -1638.4287
$ws.Range("N85").Value = -32808
$ws.Range("H98").Value = 3465.7827
$ws.Range("I98").Value = 3319.7144
$ws.Range("J98").Value = 4999.5
$ws.Range("K98").Value = 3319.7144
$ws.Range("L98").Value = 4999.5
$ws.Range("M98").Value = -1821.7144
$ws.Range("N98").Value = -7995.5
$ws.Range("H122").Value = 3465.7827
$ws.Range("I122").Value = 3319.7144
$ws.Range("J122").Value = 4999.5
$ws.Range("K122").Value = 9959.143199999999
$ws.Range("L122").Value = 14998.5
$ws.Range("M122").Value = -7509.143199999999
$ws.Range("N122").Value = -19898.5
$ws.Range("H132").Value = 14735.871
$ws.Range("I132").Value = 9460.24
$ws.Range("K132").Value = 28380.72
$ws.Range("M132").Value = -25850.72
$ws.Range("H135").Value = 943.2857
$ws.Range("I135").Value = 928.46155
$ws.Range("K135").Value = 8356.15395
$ws.Range("M135").Value = -5821.15395
$ws.Range("H138").Value = 3772.7659
$ws.Range("J138").Value = 3809.4358
$ws.Range("L138").Value = 11428.3074
$ws.Range("N138").Value = -21708.3074
$ws.Range("H141").Value = 4570
$ws.Range("I141").Value = 4665.5
$ws.Range("J141").Value = 4379
$ws.Range("K141").Value = 13996.5
$ws.Range("L141").Value = 13137
$ws.Range("M141").Value = -8816.5
$ws.Range("N141").Value = -23497

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1758112.2
$ws.Range("I61").Value = 2444.9285
$ws.Range("K61").Value = 2444.9285
$ws.Range("M61").Value = -2232.9285
$ws.Range("H132").Value = 4732899.5
$ws.Range("I132").Value = 2593
$ws.Range("K132").Value = 7779
$ws.Range("M132").Value = -5249
$ws.Range("H136").Value = 1758112.2
$ws.Range("I136").Value = 2444.9285
$ws.Range("K136").Value = 7334.7855
$ws.Range("M136").Value = -4784.7855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 471.5
$ws.Range("I22").Value = 324.57144
$ws.Range("K22").Value = 324.57144
$ws.Range("M22").Value = -151.57144
$ws.Range("H64").Value = 1010.087
$ws.Range("I64").Value = 777.5
$ws.Range("J64").Value = 1059.0526
$ws.Range("K64").Value = 777.5
$ws.Range("L64").Value = 1059.0526
$ws.Range("M64").Value = -552.5
$ws.Range("N64").Value = -1509.0526
$ws.Range("H67").Value = 1010.087
$ws.Range("I67").Value = 777.5
$ws.Range("J67").Value = 1059.0526
$ws.Range("K67").Value = 777.5
$ws.Range("L67").Value = 1059.0526
$ws.Range("M67").Value = 2.5
$ws.Range("N67").Value = -2619.0526
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 85600.94
$ws.Range("I134").Value = 135948.38
$ws.Range("J134").Value = 35253.5
$ws.Range("K134").Value = 407845.14
$ws.Range("L134").Value = 105760.5
$ws.Range("M134").Value = -405310.14
$ws.Range("N134").Value = -110830.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1998.4286
$ws.Range("I22").Value = 2831
$ws.Range("J22").Value = 1771.3636
$ws.Range("K22").Value = 2831
$ws.Range("L22").Value = 1771.3636
$ws.Range("M22").Value = -2481
$ws.Range("N22").Value = -2471.3636
$ws.Range("H58").Value = 11123.462
$ws.Range("I58").Value = 4105.7095
$ws.Range("J58").Value = 38317.25
$ws.Range("K58").Value = 4105.7095
$ws.Range("L58").Value = 38317.25
$ws.Range("M58").Value = -3902.7095
$ws.Range("N58").Value = -38723.25
$ws.Range("H122").Value = 1759.5
$ws.Range("I122").Value = 1823.8889
$ws.Range("K122").Value = 5471.6667
$ws.Range("M122").Value = -3021.6667
$ws.Range("H132").Value = 514705900
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H133").Value = 87500
$ws.Range("J133").Value = 87500
$ws.Range("L133").Value = 87500
$ws.Range("N133").Value = -92560
$ws.Range("H136").Value = 11123.462
$ws.Range("I136").Value = 4105.7095
$ws.Range("J136").Value = 38317.25
$ws.Range("K136").Value = 12317.1285
$ws.Range("L136").Value = 114951.75
$ws.Range("M136").Value = -9767.128499999999
$ws.Range("N136").Value = -120051.75
$ws.Range("H137").Value = 54500
$ws.Range("J137").Value = 54500
$ws.Range("L137").Value = 54500
$ws.Range("N137").Value = -64700
$ws.Range("H138").Value = 72333
$ws.Range("J138").Value = 72333
$ws.Range("L138").Value = 72333
$ws.Range("N138").Value = -82613

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 132.66667
$ws.Range("I2").Value = 172.14285
$ws.Range("J2").Value = 53.714287
$ws.Range("K2").Value = 1032.8571
$ws.Range("L2").Value = 322.285722
$ws.Range("M2").Value = -919.8571000000002
$ws.Range("N2").Value = -548.285722
$ws.Range("H17").Value = 2333.3333
$ws.Range("J17").Value = 2333.3333
$ws.Range("L17").Value = 6999.999899999999
$ws.Range("N17").Value = -7337.999899999999
$ws.Range("H138").Value = 2688.75
$ws.Range("I138").Value = 2688.75
$ws.Range("K138").Value = 8066.25
$ws.Range("M138").Value = -2926.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 7388.6
$ws.Range("J98").Value = 7388.6
$ws.Range("L98").Value = 7388.6
$ws.Range("N98").Value = -13378.6
$ws.Range("H122").Value = 2862.6875
$ws.Range("I122").Value = 2903.5334
$ws.Range("K122").Value = 8710.600199999999
$ws.Range("M122").Value = -6260.600199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 29139
$ws.Range("I26").Value = 28497.25
$ws.Range("K26").Value = 28497.25
$ws.Range("M26").Value = -28202.25
$ws.Range("H40").Value = 4989.7144
$ws.Range("J40").Value = 5001.3335
$ws.Range("L40").Value = 5001.3335
$ws.Range("N40").Value = -5273.3335
$ws.Range("H68").Value = 10580.637
$ws.Range("J68").Value = 1538.8
$ws.Range("L68").Value = 1538.8
$ws.Range("N68").Value = -3036.8
$ws.Range("H71").Value = 10580.637
$ws.Range("J71").Value = 1538.8
$ws.Range("L71").Value = 7694
$ws.Range("N71").Value = -15182
$ws.Range("H122").Value = 15889.75
$ws.Range("I122").Value = 19334.666
$ws.Range("J122").Value = 5555
$ws.Range("K122").Value = 58003.99800000001
$ws.Range("L122").Value = 16665
$ws.Range("M122").Value = -55553.99800000001
$ws.Range("N122").Value = -21565
$ws.Range("H132").Value = 5376815.5
$ws.Range("I132").Value = 10482
$ws.Range("J132").Value = 9976530
$ws.Range("K132").Value = 31446
$ws.Range("L132").Value = 29929590
$ws.Range("M132").Value = -28916
$ws.Range("N132").Value = -29934650

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 122810.6
$ws.Range("J62").Value = 184673
$ws.Range("L62").Value = 184673
$ws.Range("N62").Value = -185921
$ws.Range("H65").Value = 122810.6
$ws.Range("J65").Value = 184673
$ws.Range("L65").Value = 923365
$ws.Range("N65").Value = -929605
$ws.Range("H122").Value = 2749.75
$ws.Range("I122").Value = 1999.6666
$ws.Range("K122").Value = 5998.9998
$ws.Range("M122").Value = -3548.9998
$ws.Range("H132").Value = 1209166.5
$ws.Range("I132").Value = 1616.75
$ws.Range("J132").Value = 10869565
$ws.Range("K132").Value = 4850.25
$ws.Range("L132").Value = 32608695
$ws.Range("M132").Value = -2320.25
$ws.Range("N132").Value = -32613755
